$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1971.1111
$ws.Range("I12").Value = 2348.7144
$ws.Range("K12").Value = 2348.7144
$ws.Range("M12").Value = -2178.7144
$ws.Range("H33").Value = 136.42857
$ws.Range("I33").Value = 136.42857
$ws.Range("K33").Value = 136.42857
$ws.Range("M33").Value = 92.57142999999999
$ws.Range("H41").Value = 612.5
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 783.3333
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 783.3333
$ws.Range("M41").Value = 340
$ws.Range("N41").Value = -1663.3333
$ws.Range("H107").Value = 572.44446
$ws.Range("I107").Value = 521.2857
$ws.Range("K107").Value = 521.2857
$ws.Range("M107").Value = 1398.7143
$ws.Range("H113").Value = 90912950
$ws.Range("I113").Value = 125002560
$ws.Range("J113").Value = 7333.3335
$ws.Range("K113").Value = 125002560
$ws.Range("L113").Value = 7333.3335
$ws.Range("M113").Value = -124999306
$ws.Range("N113").Value = -13841.3335
$ws.Range("H125").Value = 2172.0833
$ws.Range("I125").Value = 1705.8
$ws.Range("J125").Value = 2505.1428
$ws.Range("K125").Value = 15352.2
$ws.Range("L125").Value = 22546.2852
$ws.Range("M125").Value = -12892.2
$ws.Range("N125").Value = -27466.2852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3579
$ws.Range("I45").Value = 3462.5
$ws.Range("J45").Value = 3734.3333
$ws.Range("K45").Value = 3462.5
$ws.Range("L45").Value = 3734.3333
$ws.Range("M45").Value = -3085.5
$ws.Range("N45").Value = -4488.3333
$ws.Range("H102").Value = 1262.6
$ws.Range("I102").Value = 1119.9166
$ws.Range("K102").Value = 1119.9166
$ws.Range("M102").Value = 502.0834
$ws.Range("H122").Value = 2171.9656
$ws.Range("I122").Value = 2242.9565
$ws.Range("J122").Value = 1899.8334
$ws.Range("K122").Value = 6728.869499999999
$ws.Range("L122").Value = 5699.5002
$ws.Range("M122").Value = -4278.869499999999
$ws.Range("N122").Value = -10599.5002
$ws.Range("H138").Value = 32089.375
$ws.Range("J138").Value = 32089.375
$ws.Range("L138").Value = 32089.375
$ws.Range("N138").Value = -42369.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1776.2222
$ws.Range("I99").Value = 1415.0834
$ws.Range("J99").Value = 2498.5
$ws.Range("K99").Value = 1415.0834
$ws.Range("L99").Value = 2498.5
$ws.Range("M99").Value = 82.91660000000002
$ws.Range("N99").Value = -5494.5
$ws.Range("H107").Value = 1388
$ws.Range("I107").Value = 479
$ws.Range("J107").Value = 2751.5
$ws.Range("K107").Value = 479
$ws.Range("L107").Value = 2751.5
$ws.Range("M107").Value = 1441
$ws.Range("N107").Value = -6591.5
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18745.555
$ws.Range("I31").Value = 22822.857
$ws.Range("J31").Value = 4475
$ws.Range("K31").Value = 22822.857
$ws.Range("L31").Value = 4475
$ws.Range("M31").Value = -22527.857
$ws.Range("N31").Value = -5065
$ws.Range("H34").Value = 18745.555
$ws.Range("I34").Value = 22822.857
$ws.Range("J34").Value = 4475
$ws.Range("K34").Value = 22822.857
$ws.Range("L34").Value = 4475
$ws.Range("M34").Value = -22620.857
$ws.Range("N34").Value = -4879
$ws.Range("H86").Value = 6417554
$ws.Range("I86").Value = 1497.7142
$ws.Range("J86").Value = 13902953
$ws.Range("K86").Value = 1497.7142
$ws.Range("L86").Value = 13902953
$ws.Range("M86").Value = -374.7141999999999
$ws.Range("N86").Value = -13905199
$ws.Range("H89").Value = 6417554
$ws.Range("I89").Value = 1497.7142
$ws.Range("J89").Value = 13902953
$ws.Range("K89").Value = 7488.571
$ws.Range("L89").Value = 69514765
$ws.Range("M89").Value = -1872.571
$ws.Range("N89").Value = -69525997
$ws.Range("H99").Value = 17160456
$ws.Range("I99").Value = 4169536.8
$ws.Range("J99").Value = 35718910
$ws.Range("K99").Value = 4169536.8
$ws.Range("L99").Value = 35718910
$ws.Range("M99").Value = -4168038.8
$ws.Range("N99").Value = -35721906
$ws.Range("H126").Value = 17160456
$ws.Range("I126").Value = 4169536.8
$ws.Range("J126").Value = 35718910
$ws.Range("K126").Value = 12508610.4
$ws.Range("L126").Value = 107156730
$ws.Range("M126").Value = -12506140.4
$ws.Range("N126").Value = -107161670
$ws.Range("H134").Value = 1255.5555
$ws.Range("I134").Value = 980.44446
$ws.Range("J134").Value = 1530.6666
$ws.Range("K134").Value = 2941.33338
$ws.Range("L134").Value = 4591.9998
$ws.Range("M134").Value = -406.33338
$ws.Range("N134").Value = -9661.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 599.75
$ws.Range("J17").Value = 1499
$ws.Range("L17").Value = 4497
$ws.Range("N17").Value = -4835
$ws.Range("H36").Value = 2626.8572
$ws.Range("I36").Value = 2445.75
$ws.Range("K36").Value = 7337.25
$ws.Range("M36").Value = -7168.25
$ws.Range("H39").Value = 1950
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588
$ws.Range("H131").Value = 752.64
$ws.Range("I131").Value = 512.1429000000001
$ws.Range("J131").Value = 770.74194
$ws.Range("K131").Value = 1536.4287
$ws.Range("L131").Value = 2312.22582
$ws.Range("M131").Value = 3503.5713
$ws.Range("N131").Value = -12392.22582
$ws.Range("H139").Value = 2134.5
$ws.Range("J139").Value = 3429.1667
$ws.Range("L139").Value = 10287.5001
$ws.Range("N139").Value = -20567.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2130.6365
$ws.Range("I102").Value = 2143.7
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2143.7
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -521.6999999999998
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 65371.625
$ws.Range("I132").Value = 55626.95
$ws.Range("J132").Value = 102401.4
$ws.Range("K132").Value = 166880.85
$ws.Range("L132").Value = 307204.2
$ws.Range("M132").Value = -164350.85
$ws.Range("N132").Value = -312264.2
$ws.Range("H141").Value = 52429
$ws.Range("J141").Value = 52429
$ws.Range("L141").Value = 52429
$ws.Range("N141").Value = -62789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5527.4
$ws.Range("I7").Value = 5695.4614
$ws.Range("J7").Value = 4435
$ws.Range("K7").Value = 5695.4614
$ws.Range("L7").Value = 4435
$ws.Range("M7").Value = -5583.4614
$ws.Range("N7").Value = -4659
$ws.Range("H22").Value = 5350.3335
$ws.Range("I22").Value = 5300.5
$ws.Range("J22").Value = 5375.25
$ws.Range("K22").Value = 5300.5
$ws.Range("L22").Value = 5375.25
$ws.Range("M22").Value = -5005.5
$ws.Range("N22").Value = -5965.25
$ws.Range("H27").Value = 5350.3335
$ws.Range("I27").Value = 5300.5
$ws.Range("J27").Value = 5375.25
$ws.Range("K27").Value = 5300.5
$ws.Range("L27").Value = 5375.25
$ws.Range("M27").Value = -5193.5
$ws.Range("N27").Value = -5589.25
$ws.Range("H40").Value = 4242.6924
$ws.Range("I40").Value = 4316.6665
$ws.Range("K40").Value = 4316.6665
$ws.Range("M40").Value = -4180.6665
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("H122").Value = 1964513.5
$ws.Range("I122").Value = 2804062.2
$ws.Range("J122").Value = 5566.6665
$ws.Range("K122").Value = 8412186.600000001
$ws.Range("L122").Value = 16699.9995
$ws.Range("M122").Value = -8409736.600000001
$ws.Range("N122").Value = -21599.9995
$ws.Range("H126").Value = 5527.4
$ws.Range("I126").Value = 5695.4614
$ws.Range("J126").Value = 4435
$ws.Range("K126").Value = 17086.3842
$ws.Range("L126").Value = 13305
$ws.Range("M126").Value = -14616.3842
$ws.Range("N126").Value = -18245
$ws.Range("H132").Value = 3662.625
$ws.Range("I132").Value = 4251.5
$ws.Range("J132").Value = 3466.3333
$ws.Range("K132").Value = 12754.5
$ws.Range("L132").Value = 10398.9999
$ws.Range("M132").Value = -10224.5
$ws.Range("N132").Value = -15458.9999
$ws.Range("H140").Value = 49067.5
$ws.Range("J140").Value = 49067.5
$ws.Range("L140").Value = 49067.5
$ws.Range("N140").Value = -59427.5
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4915.8335
$ws.Range("J62").Value = 5299.2
$ws.Range("L62").Value = 5299.2
$ws.Range("N62").Value = -6547.2
$ws.Range("H65").Value = 4915.8335
$ws.Range("J65").Value = 5299.2
$ws.Range("L65").Value = 26496
$ws.Range("N65").Value = -32736
$ws.Range("H81").Value = 1205.1666
$ws.Range("I81").Value = 1282.909
$ws.Range("J81").Value = 350
$ws.Range("K81").Value = 2565.818
$ws.Range("L81").Value = 700
$ws.Range("M81").Value = -1504.818
$ws.Range("N81").Value = -2822
$ws.Range("H84").Value = 1205.1666
$ws.Range("I84").Value = 1282.909
$ws.Range("J84").Value = 350
$ws.Range("K84").Value = 12829.09
$ws.Range("L84").Value = 3500
$ws.Range("M84").Value = -7525.09
$ws.Range("N84").Value = -14108
$ws.Range("H96").Value = 1137.75
$ws.Range("I96").Value = 1137.75
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1137.75
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 235.25
$ws.Range("H98").Value = 41666.668
$ws.Range("J98").Value = 41666.668
$ws.Range("L98").Value = 41666.668
$ws.Range("N98").Value = -47656.668
$ws.Range("H107").Value = 1567893.5
$ws.Range("I107").Value = 488.61905
$ws.Range("J107").Value = 5682331
$ws.Range("K107").Value = 1465.85715
$ws.Range("L107").Value = 17046993
$ws.Range("M107").Value = 454.14285
$ws.Range("N107").Value = -17050833
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H122").Value = 1835.6818
$ws.Range("I122").Value = 1870.3158
$ws.Range("K122").Value = 5610.9474
$ws.Range("M122").Value = -3160.9474
$ws.Range("H136").Value = 27028826
$ws.Range("I136").Value = 50001796
$ws.Range("J136").Value = 1800.0588
$ws.Range("K136").Value = 150005388
$ws.Range("L136").Value = 5400.1764
$ws.Range("M136").Value = -150002838
$ws.Range("N136").Value = -10500.1764
$ws.Range("N96").ClearContents()
$ws.Range("N110").ClearContents()
